# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates column G (header "K") for rows 2-36 with the actual strikeout
# counts (K) per game, replacing the previous "Strike#" derived values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 5
    3  = 5
    4  = 8
    5  = 5
    6  = 8
    7  = 5
    8  = 7
    9  = 6
    10 = 3
    11 = 7
    12 = 9
    13 = 2
    14 = 7
    15 = 4
    16 = 5
    17 = 6
    18 = 7
    19 = 6
    20 = 2
    21 = 6
    22 = 5
    23 = 5
    24 = 6
    25 = 3
    26 = 8
    27 = 2
    28 = 5
    29 = 3
    30 = 7
    31 = 5
    32 = 2
    33 = 3
    34 = 8
    35 = 5
    36 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
